$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target table (rows 2-9, columns A-T) per "Natmi following Dr Hou advice".
# Columns A-D are categorical labels (cluster / gene symbols); E-T are numeric.

$data = @(
    @("ECs",  "Ucn2", "Crhr2", "ECs", 1, 0.3333333333333333, 0.2125036666666667, 0.637511,        0.1064749515943994, 0.1064749515943994, 3, 1, 0.7693313333333333, 2.307994,        0.2551791101080821, 0.2551791101080821, 0.1634857292148889, 1.471371562934,  0.02717018339665996, 0.02717018339665996),
    @("ECs",  "Ucn2", "Crhr2", "sCs", 1, 0.3333333333333333, 0.2125036666666667, 0.637511,        0.1064749515943994, 0.1064749515943994, 3, 1, 2.245536666666666,  6.736609999999999, 0.7448208898919179, 0.7448208898919177, 0.477184775301111,  4.29466297771,  0.07930476819773943, 0.07930476819773943),
    @("FAPs", "Ucn2", "Crhr2", "ECs", 3, 1,                  1.193141333333333, 3.579424,         0.597823405613129,  0.597823405613129,  3, 1, 0.7693313333333333, 2.307994,        0.2551791101080821, 0.2551791101080821, 0.9179210128284444, 8.261289115456,  0.1525520446461413,  0.1525520446461413),
    @("FAPs", "Ucn2", "Crhr2", "sCs", 3, 1,                  1.193141333333333, 3.579424,         0.597823405613129,  0.597823405613129,  3, 1, 2.245536666666666,  6.736609999999999, 0.7448208898919179, 0.7448208898919177, 2.679242612515555,  24.11318351264, 0.4452713609669877,  0.4452713609669876),
    @("M2",   "Ucn2", "Crhr2", "ECs", 3, 1,                  0.1501633333333333, 0.45049,         0.07523933068411523, 0.07523933068411523, 3, 1, 0.7693313333333333, 2.307994,        0.2551791101080821, 0.2551791101080821, 0.1155253574511111, 1.03972821706, 0.01919950544910024, 0.01919950544910024),
    @("M2",   "Ucn2", "Crhr2", "sCs", 3, 1,                  0.1501633333333333, 0.45049,         0.07523933068411523, 0.07523933068411523, 3, 1, 2.245536666666666,  6.736609999999999, 0.7448208898919179, 0.7448208898919177, 0.3371972709888889, 3.0347754389,  0.05603982523501499, 0.05603982523501498),
    @("sCs",  "Ucn2", "Crhr2", "ECs", 3, 1,                  0.4400006666666667, 1.320002,        0.2204623121083564, 0.2204623121083564, 3, 1, 0.7693313333333333, 2.307994,        0.2551791101080821, 0.2551791101080821, 0.3385062995542222, 3.046556695988, 0.05625737661618065, 0.05625737661618065),
    @("sCs",  "Ucn2", "Crhr2", "sCs", 3, 1,                  0.4400006666666667, 1.320002,        0.2204623121083564, 0.2204623121083564, 3, 1, 2.245536666666666,  6.736609999999999, 0.7448208898919179, 0.7448208898919177, 0.9880376303577777, 8.892338673219999, 0.1642049354921757, 0.1642049354921757)
)

$rowIndex = 2
foreach ($row in $data) {
    $colIndex = 1
    foreach ($val in $row) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex++
    }
    $rowIndex++
}
